$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(4)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Paragraph 2: split the trailing run into three runs ---
# Original run 5 text: " (and LaTeX) run in a Windows computer circa 1998."
# becomes:
#   run 5: " (and "                                    (plain)
#   run 6: "LaTeX"                                      (sz=1200, Consolas)
#   run 7: ") run in a Windows computer circa 1998."    (plain)
$para2 = $tr.Paragraphs(2, 1)
$run5 = $para2.Runs(5, 1)
$run5.Text = " (and "

# Insert the two new runs first (they inherit plain formatting from run5),
# then apply the Consolas/12pt styling only to the "LaTeX" run afterwards,
# so the trailing run stays un-styled.
$run5again = $para2.Runs(5, 1)
$run5again.InsertAfter("LaTeX") | Out-Null

$para2b = $tr.Paragraphs(2, 1)
$run6 = $para2b.Runs(6, 1)
$run6.InsertAfter(") run in a Windows computer circa 1998.") | Out-Null

$para2c = $tr.Paragraphs(2, 1)
$run6again = $para2c.Runs(6, 1)
$run6again.Font.Size = 12
$run6again.Font.Name = "Consolas"

# --- Paragraph 3: replace the whole sentence ---
$para3 = $tr.Paragraphs(3, 1)
$run3 = $para3.Runs(1, 1)
$run3.Text = "Jupyter implements the literate programming paradigm, but I haven" + [char]0x2019 + "t seen the Markdown part gain a lot of traction."
